$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.737.68"

$ws.Range("D3").Value = "3.629.60"
$ws.Range("E3").Value = "  +6.17%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'592.27"
$ws.Range("E5").Value = "  +4.41%  "

$ws.Range("D6").Value = "'191.82"
$ws.Range("E6").Value = "  +7.64%  "

$ws.Range("D7").Value = "'0.652"
$ws.Range("E7").Value = "  +3.03%  "

$ws.Range("D8").Value = "3.618.99"
$ws.Range("E8").Value = "  +6.16%  "

$ws.Range("E9").Value = "  +0.03%  "

$ws.Range("E10").Value = "  +2.78%  "

$ws.Range("D11").Value = "'0.665"
$ws.Range("E11").Value = "  +4.04%  "

$ws.Range("D12").Value = "'58.16"
$ws.Range("E12").Value = "  +6.56%  "

$ws.Range("D13").Value = "'0.0000297"
$ws.Range("E13").Value = "  +6.14%  "

$ws.Range("D14").Value = "'9.79"
$ws.Range("E14").Value = "  +5.02%  "

$ws.Range("D15").Value = "4.212.97"
$ws.Range("E15").Value = "  +6.62%  "

$ws.Range("D16").Value = "3.640.25"
$ws.Range("E16").Value = "  +6.51%  "

$ws.Range("E17").Value = "  +5.91%  "

$ws.Range("D18").Value = "70.691.89"
$ws.Range("E18").Value = "  +6.85%  "

$ws.Range("D19").Value = "'12.59"
$ws.Range("E19").Value = "  +5.14%  "

$ws.Range("E20").Value = "  +0.66%  "

$ws.Range("D21").Value = "'1.05"
$ws.Range("E21").Value = "  +4.76%  "

$ws.Range("D22").Value = "'494.46"
$ws.Range("E22").Value = "  +5.96%  "

$ws.Range("D23").Value = "'5.40"
$ws.Range("E23").Value = "  +9.20%  "

$ws.Range("D24").Value = "'17.25"
$ws.Range("E24").Value = "  +17.03%  "

$ws.Range("D25").Value = "'4.50"
$ws.Range("E25").Value = "  +8.68%  "

$ws.Range("D26").Value = "'91.06"
$ws.Range("E26").Value = "  +1.15%  "

$ws.Range("E27").Value = "  +6.30%  "

$ws.Range("D28").Value = "'11.27"
$ws.Range("E28").Value = "  +4.72%  "

$ws.Range("D29").Value = "'9.46"
$ws.Range("E29").Value = "  +7.05%  "

$ws.Range("D30").Value = "'32.39"
$ws.Range("E30").Value = "  +2.80%  "

$ws.Range("E31").Value = "  +12.05%  "

$ws.Range("D32").Value = "'12.23"
$ws.Range("E32").Value = "  +5.76%  "

$ws.Range("D33").Value = "'620.30"
$ws.Range("E33").Value = "  +6.31%  "

$ws.Range("E34").Value = "  +8.04%  "

$ws.Range("D35").Value = "'65.21"
$ws.Range("E35").Value = "  +4.28%  "

$ws.Range("B36").Value = "PEPE"
$ws.Range("C36").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D36").Value = "0.0₃0829"
$ws.Range("E36").Value = "  +8.09%  "

$ws.Range("B37").Value = "TheGraph"
$ws.Range("C37").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D37").Value = "'0.414"
$ws.Range("E37").Value = "  +8.27%  "

$ws.Range("B38").Value = "InjectiveProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D38").Value = "'38.30"
$ws.Range("E38").Value = "  +5.03%  "

$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "'0.148"
$ws.Range("E39").Value = "  +2.24%  "

$ws.Range("D40").Value = "'0.999"
$ws.Range("E40").Value = "  -0.10%  "

$ws.Range("E41").Value = "  +2.78%  "

$ws.Range("D42").Value = "3.333.54"
$ws.Range("E42").Value = "  +6.51%  "

$ws.Range("D43").Value = "'3.07"
$ws.Range("E43").Value = "  +4.86%  "

$ws.Range("D44").Value = "'0.0449"
$ws.Range("E44").Value = "  +6.18%  "

$ws.Range("E45").Value = "  +8.24%  "

$ws.Range("D46").Value = "'3.36"
$ws.Range("E46").Value = "  +5.90%  "

$ws.Range("D47").Value = "'0.138"
$ws.Range("E47").Value = "  +2.86%  "

$ws.Range("D48").Value = "'9.18"
$ws.Range("E48").Value = "  +7.09%  "

$ws.Range("D49").Value = "'2.73"
$ws.Range("E49").Value = "  +3.59%  "

$ws.Range("D50").Value = "'3.32"
$ws.Range("E50").Value = "  +4.67%  "

$ws.Range("D51").Value = "'1.00"
$ws.Range("E51").Value = "  -0.03%  "
